$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "depth" attribute row (row 9) to "depth_mat_file"
$ws.Range("A9").Value = "depth_mat_file"

# Insert a new row after it for the "depth_API_bottle_summary" attribute,
# shifting the remaining attribute rows down by one
$ws.Rows.Item(10).Insert()
$ws.Range("A10").Value = "depth_API_bottle_summary"
$ws.Range("B10").Value = "Depth of sample below sea surface  "
$ws.Range("C10").Value = "numeric"
$ws.Range("D10").Value = "meter"
